$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translated_Sheet1")
$ws.Activate()

# Overwrite column A text with the final (untranslated-label) English values
$ws.Range("A1").Value = "Hungarian column"
$ws.Range("A2").Value = "God bless the seed."
$ws.Range("A3").Value = "It's only the Fidesz."
$ws.Range("A4").Value = "State deficit in proportion to GDP"

# Drop the now-unused column B entirely
$ws.Range("B1:B4").ClearContents()
$ws.Columns.Item(2).Delete()

# Column A sizing (bestFit width ~32 characters)
$ws.Columns.Item(1).ColumnWidth = 31.1

# Row height / sheet format
$ws.Rows.Item(1).RowHeight = 14.5
$ws.Rows.Item(2).RowHeight = 14.5
$ws.Rows.Item(3).RowHeight = 14.5
$ws.Rows.Item(4).RowHeight = 14.5
$ws.StandardHeight = 14.5

# Clear the header-row bold/border style (s="1" -> default)
$ws.Range("A1").Style = "Normal"

# Page margins (Excel COM uses points; 72pt = 1in)
$ws.PageSetup.LeftMargin = 50.4
$ws.PageSetup.RightMargin = 50.4
$ws.PageSetup.TopMargin = 54
$ws.PageSetup.BottomMargin = 54
$ws.PageSetup.HeaderMargin = 21.599999999999998
$ws.PageSetup.FooterMargin = 21.599999999999998

Write-Output "done"
